$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2928702470674125
$ws.Range("C2").Value = 0.04022481641014508
$ws.Range("D2").Value = 0.3101056510207769
$ws.Range("F2").Value = 1.227379493136226
$ws.Range("G2").Value = 0.5897553890308203
$ws.Range("H2").Value = 0.7251471904914908
$ws.Range("J2").Value = 0.3102659089815063
$ws.Range("K2").Value = 0.2565142136828911
$ws.Range("M2").Value = 0.2585865407620602
$ws.Range("N2").Value = 1.574638704693657
$ws.Range("O2").Value = 2.609909659457458
$ws.Range("B3").Value = 0.2616281610741282
$ws.Range("C3").Value = 0.03733829473171824
$ws.Range("D3").Value = 0.3046319723753612
$ws.Range("F3").Value = 1.229998089562862
$ws.Range("G3").Value = 0.592246784886818
$ws.Range("H3").Value = 0.7296350741670992
$ws.Range("J3").Value = 0.3085395197972147
$ws.Range("K3").Value = 0.225078018016788
$ws.Range("M3").Value = 0.2454924147735937
$ws.Range("N3").Value = 1.589794889120018
$ws.Range("O3").Value = 2.624353282935544
$ws.Range("B4").Value = 0.2424630465902169
$ws.Range("C4").Value = 0.03555506210097548
$ws.Range("D4").Value = 0.3014027707341427
$ws.Range("F4").Value = 1.232202206187594
$ws.Range("G4").Value = 0.5941175558829528
$ws.Range("H4").Value = 0.7326613133573545
$ws.Range("J4").Value = 0.3076347288732464
$ws.Range("K4").Value = 0.2057429455475841
$ws.Range("M4").Value = 0.2375511242076769
$ws.Range("N4").Value = 1.599591373639397
$ws.Range("O4").Value = 2.634503127322134
$ws.Range("B5").Value = 0.234658061941019
$ws.Range("C5").Value = 0.03482569102951771
$ws.Range("D5").Value = 0.3001200965987891
$ws.Range("F5").Value = 1.233250505642843
$ws.Range("G5").Value = 0.5949656819441671
$ws.Range("H5").Value = 0.7339626691662673
$ws.Range("J5").Value = 0.3073050980170322
$ws.Range("K5").Value = 0.1978559223284719
$ws.Range("M5").Value = 0.2343399884964299
$ws.Range("N5").Value = 1.603706905703778
$ws.Range("O5").Value = 2.638961678269354
$ws.Range("B6").Value = 0.2333623646482295
$ws.Range("C6").Value = 0.03470441858616624
$ws.Range("D6").Value = 0.2999091221197574
$ws.Range("F6").Value = 1.23343364566481
$ws.Range("G6").Value = 0.5951116934266309
$ws.Range("H6").Value = 0.734182876046674
$ws.Range("J6").Value = 0.3072527247303611
$ws.Range("K6").Value = 0.1965458304101304
$ws.Range("M6").Value = 0.2338082986057373
$ws.Range("N6").Value = 1.604397739984472
$ws.Range("O6").Value = 2.639721495387846
$ws.Range("B7").Value = 0.2423577650912421
$ws.Range("C7").Value = 0.03554523636988449
$ws.Range("D7").Value = 0.3013853373149544
$ws.Range("F7").Value = 1.232215735939349
$ws.Range("G7").Value = 0.5941286467038083
$ws.Range("H7").Value = 0.7326785879000894
$ws.Range("J7").Value = 0.3076301250644562
$ws.Range("K7").Value = 0.2056366094058291
$ws.Range("M7").Value = 0.2375077161712298
$ws.Range("N7").Value = 1.599646377510648
$ws.Range("O7").Value = 2.634561951236307
$ws.Range("B8").Value = 0.2820946121734664
$ws.Range("C8").Value = 0.03923183493475335
$ws.Range("D8").Value = 0.3081910854831165
$ws.Range("F8").Value = 1.228158720962256
$ws.Range("G8").Value = 0.5905436458113869
$ws.Range("H8").Value = 0.7266384844231908
$ws.Range("J8").Value = 0.3096384616568244
$ws.Range("K8").Value = 0.2456821809446552
$ws.Range("M8").Value = 0.2540513597385043
$ws.Range("N8").Value = 1.579762772711874
$ws.Range("O8").Value = 2.614623976614709
$ws.Range("B9").Value = 0.3601395995491146
$ws.Range("C9").Value = 0.04637287521481426
$ws.Range("D9").Value = 0.3225763672660662
$ws.Range("F9").Value = 1.224928133159075
$ws.Range("G9").Value = 0.586219417513405
$ws.Range("H9").Value = 0.7169380578896209
$ws.Range("J9").Value = 0.3148068691886152
$ws.Range("K9").Value = 0.3239301414819522
$ws.Range("M9").Value = 0.2872674146697989
$ws.Range("N9").Value = 1.544660931138859
$ws.Range("O9").Value = 2.585686036587859
$ws.Range("B10").Value = 0.4175332327528736
$ws.Range("C10").Value = 0.05156343437006683
$ws.Range("D10").Value = 0.3337728740830386
$ws.Range("F10").Value = 1.225428660442986
$ws.Range("G10").Value = 0.5846926671792829
$ws.Range("H10").Value = 0.7111139291148447
$ws.Range("J10").Value = 0.3193527662455722
$ws.Range("K10").Value = 0.3812280201049134
$ws.Range("M10").Value = 0.3121350182856446
$ws.Range("N10").Value = 1.521239082146819
$ws.Range("O10").Value = 2.570611902147562
$ws.Range("B11").Value = 0.4436509051968187
$ws.Range("C11").Value = 0.05391218693135613
$ws.Range("D11").Value = 0.3390015678210148
$ws.Range("F11").Value = 1.226279145521943
$ws.Range("G11").Value = 0.5843566341307138
$ws.Range("H11").Value = 0.7087463582228821
$ws.Range("J11").Value = 0.3215831715872781
$ws.Range("K11").Value = 0.4072490497424042
$ws.Range("M11").Value = 0.3235470679448866
$ws.Range("N11").Value = 1.511097130797134
$ws.Range("O11").Value = 2.565096212244384
$ws.Range("B12").Value = 0.4535418281062107
$ws.Range("C12").Value = 0.05479976107154982
$ws.Range("D12").Value = 0.3410008711612704
$ws.Range("F12").Value = 1.22669062314354
$ws.Range("G12").Value = 0.5842809404880143
$ws.Range("H12").Value = 0.7078902776471665
$ws.Range("J12").Value = 0.3224510933164737
$ws.Range("K12").Value = 0.4170957536605613
$ws.Range("M12").Value = 0.327882650941028
$ws.Range("N12").Value = 1.507330340367446
$ws.Range("O12").Value = 2.563200343009328
$ws.Range("B13").Value = 0.4514116176269738
$ws.Range("C13").Value = 0.05460868915645278
$ws.Range("D13").Value = 0.3405694290028976
$ws.Range("F13").Value = 1.226598029441803
$ws.Range("G13").Value = 0.5842949495009293
$ws.Range("H13").Value = 0.7080728511333803
$ws.Range("J13").Value = 0.3222631345962981
$ws.Range("K13").Value = 0.4149754053668744
$ws.Range("M13").Value = 0.3269482828927437
$ws.Range("N13").Value = 1.508138305985057
$ws.Range("O13").Value = 2.563600079317297
$ws.Range("B14").Value = 0.4444646260709249
$ws.Range("C14").Value = 0.05398524553696404
$ws.Range("D14").Value = 0.3391656655716702
$ws.Range("F14").Value = 1.226311206777225
$ws.Range("G14").Value = 0.5843493735147831
$ws.Range("H14").Value = 0.7086751172984265
$ws.Range("J14").Value = 0.3216541090097707
$ws.Range("K14").Value = 0.4080592853951543
$ws.Range("M14").Value = 0.3239034780245191
$ws.Range("N14").Value = 1.510785756619592
$ws.Range("O14").Value = 2.564936374570095
$ws.Range("B15").Value = 0.440209470926078
$ws.Range("C15").Value = 0.05360312606339335
$ws.Range("D15").Value = 0.3383083302102818
$ws.Range("F15").Value = 1.226147160244935
$ws.Range("G15").Value = 0.5843894238524712
$ws.Range("H15").Value = 0.7090492908594968
$ws.Range("J15").Value = 0.3212840982924945
$ws.Range("K15").Value = 0.4038220476461163
$ws.Range("M15").Value = 0.3220402745233386
$ws.Range("N15").Value = 1.512417000917337
$ws.Range("O15").Value = 2.565779999071111
$ws.Range("B16").Value = 0.4158265142455946
$ws.Range("C16").Value = 0.05140968238080745
$ws.Range("D16").Value = 0.3334338771721548
$ws.Range("F16").Value = 1.225385599519569
$ws.Range("G16").Value = 0.5847218418777373
$ws.Range("H16").Value = 0.7112743216650159
$ws.Range("J16").Value = 0.3192102691591145
$ws.Range("K16").Value = 0.3795265501803158
$ws.Range("M16").Value = 0.3113911985531601
$ws.Range("N16").Value = 1.52191220159135
$ws.Range("O16").Value = 2.570999351914367
$ws.Range("B17").Value = 0.4008702651589999
$ws.Range("C17").Value = 0.05006084491400031
$ws.Range("D17").Value = 0.3304781155024585
$ws.Range("F17").Value = 1.225077819831313
$ws.Range("G17").Value = 0.5850175912735835
$ws.Range("H17").Value = 0.71271145063983
$ws.Range("J17").Value = 0.317979618594876
$ws.Range("K17").Value = 0.3646103661802442
$ws.Range("M17").Value = 0.3048836864231532
$ws.Range("N17").Value = 1.527868538544823
$ws.Range("O17").Value = 2.574544789797216
$ws.Range("B18").Value = 0.3922687094544415
$ws.Range("C18").Value = 0.04928385993262907
$ws.Range("D18").Value = 0.3287907881593668
$ws.Range("F18").Value = 1.224959421672125
$ws.Range("G18").Value = 0.5852214414012451
$ws.Range("H18").Value = 0.7135645824540404
$ws.Range("J18").Value = 0.3172870745978713
$ws.Range("K18").Value = 0.3560268577093098
$ws.Range("M18").Value = 0.3011501318212169
$ws.Range("N18").Value = 1.531342729083365
$ws.Range("O18").Value = 2.576710317686405
$ws.Range("B19").Value = 0.3893565395357825
$ws.Range("C19").Value = 0.04902058703557088
$ws.Range("D19").Value = 0.3282216823254345
$ws.Range("F19").Value = 1.224929407843177
$ws.Range("G19").Value = 0.5852962567990332
$ws.Range("H19").Value = 0.7138579972672687
$ws.Range("J19").Value = 0.3170552191657379
$ws.Range("K19").Value = 0.3531199396820739
$ws.Range("M19").Value = 0.2998876352305047
$ws.Range("N19").Value = 1.532527319222549
$ws.Range("O19").Value = 2.577465221542354
$ws.Range("B20").Value = 0.4024622961258899
$ws.Range("C20").Value = 0.05020455239272792
$ws.Range("D20").Value = 0.3307914430479144
$ws.Range("F20").Value = 1.225104516674783
$ws.Range("G20").Value = 0.5849826159519438
$ws.Range("H20").Value = 0.7125557201745352
$ws.Range("J20").Value = 0.3181090409378271
$ws.Range("K20").Value = 0.3661986492588198
$ws.Range("M20").Value = 0.3055754516818894
$ws.Range("N20").Value = 1.527229481185563
$ws.Range("O20").Value = 2.57415430263498
$ws.Range("B21").Value = 0.4465051106984674
$ws.Range("C21").Value = 0.0541684166252594
$ws.Range("D21").Value = 0.3395774620533984
$ws.Range("F21").Value = 1.226393027899533
$ws.Range("G21").Value = 0.5843319886531191
$ws.Range("H21").Value = 0.7084971193618372
$ws.Range("J21").Value = 0.3218323621695021
$ws.Range("K21").Value = 0.4100909090732898
$ws.Range("M21").Value = 0.3247974306270862
$ws.Range("N21").Value = 1.510006134633931
$ws.Range("O21").Value = 2.564538640601398
$ws.Range("B22").Value = 0.4752936096700182
$ws.Range("C22").Value = 0.05674823662047856
$ws.Range("D22").Value = 0.3454321187944345
$ws.Range("F22").Value = 1.227756272174517
$ws.Range("G22").Value = 0.5842072681480772
$ws.Range("H22").Value = 0.7060804334306567
$ws.Range("J22").Value = 0.3244016479738292
$ws.Range("K22").Value = 0.4387365875436728
$ws.Range("M22").Value = 0.3374421045714797
$ws.Range("N22").Value = 1.499179539565858
$ws.Range("O22").Value = 2.559377996198066
$ws.Range("B23").Value = 0.459928470876207
$ws.Range("C23").Value = 0.05537234428354054
$ws.Range("D23").Value = 0.3422971362851399
$ws.Range("F23").Value = 1.226981046088653
$ws.Range("G23").Value = 0.5842463367601454
$ws.Range("H23").Value = 0.7073487052317091
$ws.Range("J23").Value = 0.323017953738173
$ws.Range("K23").Value = 0.4234517358918879
$ws.Range("M23").Value = 0.3306859814279619
$ws.Range("N23").Value = 1.553740884834326
$ws.Range("O23").Value = 2.562029545882041
$ws.Range("B24").Value = 0.4017425482554131
$ws.Range("C24").Value = 0.05013958697323062
$ws.Range("D24").Value = 0.3306497503387078
$ws.Range("F24").Value = 1.225092264647948
$ws.Range("G24").Value = 0.5849983229490618
$ws.Range("H24").Value = 0.7126260420989894
$ws.Range("J24").Value = 0.3180504824528896
$ws.Range("K24").Value = 0.3654806113422353
$ws.Range("M24").Value = 0.3052626806035263
$ws.Range("N24").Value = 1.52751824385367
$ws.Range("O24").Value = 2.574330445627368
$ws.Range("B25").Value = 0.3390154454519063
$ws.Range("C25").Value = 0.04445071153891433
$ws.Range("D25").Value = 0.3185741192665148
$ws.Range("F25").Value = 1.225296904603596
$ws.Range("G25").Value = 0.5870994774192297
$ws.Range("H25").Value = 0.7193331838681871
$ws.Range("J25").Value = 0.313277115092923
$ws.Range("K25").Value = 0.3027941149423725
$ws.Range("M25").Value = 0.2781995619600579
$ws.Range("N25").Value = 1.553740884834326
$ws.Range("O25").Value = 2.592427549177359
